$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Rows 1-12: simple value swaps
$t.Cell(1,1).Range.Text  = "0M"
$t.Cell(2,1).Range.Text  = "0M"
$t.Cell(3,1).Range.Text  = "0M"
$t.Cell(4,1).Range.Text  = "588"
$t.Cell(5,1).Range.Text  = "0.00003"
$t.Cell(6,1).Range.Text  = "0.00021"
$t.Cell(7,1).Range.Text  = "0.00006"
$t.Cell(8,1).Range.Text  = "0.00002"
$t.Cell(9,1).Range.Text  = "0.00006"
$t.Cell(10,1).Range.Text = "0.00011"
$t.Cell(11,1).Range.Text = "0.00012"
$t.Cell(12,1).Range.Text = "0.03805"

# Rows 44-46: collapse tab-separated summary rows back to single values
$t.Cell(44,1).Range.Text = "99.98"
$t.Cell(45,1).Range.Text = "0.04"
$t.Cell(46,1).Range.Text = "180"
